# The header row (row 1, columns B:Y) used to hold shared-string labels
# like "2023 год", "2022 год", ... . Replace each with the plain numeric
# year value (2000 in column B through 2023 in column Y), which also
# collapses sharedStrings.xml down to the single remaining string used
# by A2 ("Индекс потребительской уверенности").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$years = 2000..2023
for ($i = 0; $i -lt $years.Count; $i++) {
    $col = $i + 2  # column B = 2 ... column Y = 25
    $ws.Cells.Item(1, $col).Value = $years[$i]
}

# Match the author's resulting selection/scroll state: B1:Y1 selected
# with B1 as the active cell, scrolled right so column E is at the
# left edge of the viewport.
$ws.Range("B1:Y1").Select()
$excel.ActiveWindow.ScrollColumn = 5
